$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'97.351.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.67%  '
$ws.Range("D3").Value = "'3.740.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.01%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'239.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.47%  '
$ws.Range("E6").Value = '  +0.81%  '
$ws.Range("D7").Value = "'662.05"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.64%  '
$ws.Range("D8").Value = "'0.442"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.81%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = "'1.07"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.90%  '
$ws.Range("B10").Value = 'USDC'
$ws.Range("C10").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D10").Value = "'0.999"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.01%  '
$ws.Range("D11").Value = "'3.739.43"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.02%  '
$ws.Range("E12").Value = '  +18.34%  '
$ws.Range("D13").Value = "'45.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.69%  '
$ws.Range("E14").Value = '  +1.52%  '
$ws.Range("D15").Value = "'6.96"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.71%  '
$ws.Range("D16").Value = "'4.436.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.97%  '
$ws.Range("D17").Value = "'97.055.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.59%  '
$ws.Range("D18").Value = "'9.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.51%  '
$ws.Range("D19").Value = "'3.735.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.55%  '
$ws.Range("D20").Value = "'13.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.06%  '
$ws.Range("D21").Value = "'18.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.69%  '
$ws.Range("E22").Value = '  -3.89%  '
$ws.Range("D23").Value = "'530.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.47%  '
$ws.Range("D24").Value = "'3.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("D25").Value = "'0.0000229"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +11.83%  '
$ws.Range("E26").Value = '  -3.26%  '
$ws.Range("D27").Value = "'107.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.59%  '
$ws.Range("E28").Value = '  +14.59%  '
$ws.Range("E29").Value = '  +1.56%  '
$ws.Range("D30").Value = "'3.932.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.79%  '
$ws.Range("D31").Value = "'12.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.48%  '
$ws.Range("E32").Value = '  +0.64%  '
$ws.Range("E33").Value = '  +0.09%  '
$ws.Range("E34").Value = '  +3.78%  '
$ws.Range("D35").Value = "'1.85"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.89%  '
$ws.Range("D36").Value = "'32.71"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.46%  '
$ws.Range("D37").Value = "'0.996"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.43%  '
$ws.Range("D38").Value = "'649.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.35%  '
$ws.Range("D39").Value = "'0.597"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.11%  '
$ws.Range("E40").Value = '  -0.51%  '
$ws.Range("D42").Value = "'0.166"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.12%  '
$ws.Range("D43").Value = "'6.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.90%  '
$ws.Range("B44").Value = 'ImmutableX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D44").Value = "'2.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.82%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = "'40.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.40%  '
$ws.Range("D46").Value = "'0.981"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.30%  '
$ws.Range("D47").Value = "'0.480"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +10.26%  '
$ws.Range("D48").Value = "'0.0461"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("D49").Value = "'2.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.82%  '
$ws.Range("D50").Value = "'23.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'8.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.46%  '
